# Append the 2025-09-19 tracker snapshot (rows 62-66) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$goals = @(
    @{ Id = "G2"; Name = "Workout" },
    @{ Id = "G3"; Name = "Eat Healthy" },
    @{ Id = "G4"; Name = "Read Book" },
    @{ Id = "G5"; Name = "Investment Plan" },
    @{ Id = "G6"; Name = "Spend 10 Hours without phone" }
)

$startRow = 62
$date = 45919
$progress = 0.8874492252651535
$percentage = 0
$change = -0.01

for ($i = 0; $i -lt $goals.Count; $i++) {
    $row = $startRow + $i
    $goal = $goals[$i]

    $ws.Cells.Item($row, 1).Value = $goal.Id
    $ws.Cells.Item($row, 2).Value = $goal.Name
    $ws.Cells.Item($row, 3).Value = $date
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 4).Value = $progress
    $ws.Cells.Item($row, 5).Value = $percentage
    $ws.Cells.Item($row, 6).Value = $change
}
